$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new server data row (row 2), writing cells in the order that
# reproduces the author's shared-string insertion order: IP, ID/Name, ServerID.
$ws.Range("F2").Value = "127.0.0.1"

$ws.Range("A2").Value = "GameServer_1"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "GameServer_1"

$ws.Range("B2").Value = "000104001"

$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1

$ws.Range("F2").NumberFormat = "@"

$ws.Range("G2").Value = 4001

# Move the data validation dropdown down so it no longer covers the newly
# filled-in F2 cell.
$ws.Range("F2:F1048576").Validation.Delete() | Out-Null
$ws.Range("F3:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"') | Out-Null

# Update the active selection to match the saved view state.
$ws.Range("G3").Select() | Out-Null
